$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.022.30'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '3.263.76'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.10'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.47'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -3.65%  '
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.407'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.51%  '
$ws.Range("D12").Value = '3.833.75'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = '68.029.33'
$ws.Range("E14").Value = '  -1.10%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.33'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.03%  '
$ws.Range("E16").Value = '  -2.82%  '
$ws.Range("D17").Value = '3.285.28'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("E18").Value = '  -2.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.25'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '415.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +4.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.52'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.19'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.507'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("E25").Value = '  -3.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.186'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.36'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.91%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.95'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.56'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("E31").Value = '  -4.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.83'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.82%  '
$ws.Range("E33").Value = '  -5.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '163.36'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.43'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.35%  '
$ws.Range("E36").Value = '  -4.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.79'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.794'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.43'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.21%  '
$ws.Range("E40").Value = '  -5.12%  '
$ws.Range("D41").Value = '2.639.06'
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.42'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0672'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '337.14'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.16'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.18%  '
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.22'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.974'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("E49").Value = '  -1.97%  '
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '30.56'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.42%  '
